# Generate Report for Handoff
# Updates the localization-status workbook to reflect that "b.md" has been
# handed off again (new xliff files generated), superseding the previous
# "Handed back: in sync with en-US" status.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-04 16:42:26"

# --- zh-cn detail sheet -------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-04 16:42:22"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f826c960ad716232de918ed72c3568cf89b121b5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75573dc9371deae830e6947cf1e2dcda87ad23df/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de detail sheet -------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-04 16:42:26"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f826c960ad716232de918ed72c3568cf89b121b5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75573dc9371deae830e6947cf1e2dcda87ad23df/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
